# "separate dept from affiliations"
#
# 1. Sheet "dept hours" -> new sheet "unit(accumulative) hours" (copy of the
#    original dept/affiliation accumulative breakdown, with header relabeled).
# 2. Sheet "dept hours" itself is renamed "department hours" and recomputed
#    to show hours/percentage per single primary department (no more CSL/AE
#    accumulative rows).
# 3. Sheet "PI hours" keeps a simplified single "dept" per PI and gains a new
#    "app" column holding the full affiliation list that used to live in dept.

$wb = $excel.ActiveWorkbook

$piSheet = $wb.Worksheets.Item(1)
$deptSheet = $wb.Worksheets.Item(2)

# --- Step 1: duplicate the current "dept hours" sheet to the end of the
# workbook; it becomes the new accumulative-unit breakdown sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$deptSheet.Copy($null, $lastSheet)
$unitSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$unitSheet.Name = "unit(accumulative) hours"
$unitSheet.Range("B1").Value = "unit(accumulative)"

# --- Step 2: rename & rebuild the "dept hours" sheet as "department hours"
# with per-PI-primary-department totals.
$deptSheet.Name = "department hours"

$deptSheet.Range("B1").Value = "dept"
$deptSheet.Range("C1").Value = "hours"
$deptSheet.Range("D1").Value = "percentage"

$deptSheet.Range("A2").Value = 0
$deptSheet.Range("B2").Value = "ABE"
$deptSheet.Range("C2").Value = 47.5
$deptSheet.Range("D2").Value = 45.45454545454545

$deptSheet.Range("A3").Value = 1
$deptSheet.Range("B3").Value = "ECE"
$deptSheet.Range("C3").Value = 43
$deptSheet.Range("D3").Value = 41.14832535885167

$deptSheet.Range("A4").Value = 2
$deptSheet.Range("B4").Value = "ME"
$deptSheet.Range("C4").Value = 14
$deptSheet.Range("D4").Value = 13.39712918660287

# the old sheet had 5 data rows (rows 2-6); the new data only needs 3
# (rows 2-4), so drop the now-unused trailing rows.
$deptSheet.Rows.Item(5).Delete()
$deptSheet.Rows.Item(5).Delete()

# --- Step 3: on "PI hours", add the new "app" column (the old affiliation
# list) and simplify "dept" down to each PI's single primary department.

$piSheet.Range("F1").Value = "app"
$appHeader = $piSheet.Range("F1")
$appHeader.Font.Bold = $true
$appHeader.HorizontalAlignment = -4108
$appHeader.VerticalAlignment = -4160
$appHeader.Borders.LineStyle = 1

$piSheet.Range("E2").Value = "ABE"
$piSheet.Range("F2").Value = "['ABE', 'CSL']"

$piSheet.Range("E3").Value = "ECE"
$piSheet.Range("F3").Value = "['ECE', 'CSL']"

$piSheet.Range("E4").Value = "ME"
$piSheet.Range("F4").Value = "['ME', 'AE', 'CSL']"

$piSheet.Range("E5").Value = "ECE"
$piSheet.Range("F5").Value = "['ECE', 'CSL']"
